$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new column (item_num) before the existing last column
#    (column O, "comment") on every data sheet. This pushes the old
#    "comment" column from O to P and creates a fresh, empty O column.
# ------------------------------------------------------------------
$sheetNames = @("Measures", "ID", "Dems", "Dates", "NewVars")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(15).Insert()
    $ws.Cells.Item(1, 15).Value = "item_num"
}

# ------------------------------------------------------------------
# 2. Fill in item_num = 1 for every data row on the Measures sheet.
# ------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("Measures")
for ($r = 2; $r -le 76; $r++) {
    $wsMeasures.Cells.Item($r, 15).Value = 1
}

# ------------------------------------------------------------------
# 3. Update the _FilterDatabase defined names so their range keeps up
#    with the extra column (O -> P).
# ------------------------------------------------------------------
foreach ($dn in $wb.Names) {
    if ($dn.Name -eq "Dems!_FilterDatabase") {
        $dn.RefersTo = "=Dems!`$A`$1:`$P`$1"
    }
    elseif ($dn.Name -eq "ID!_FilterDatabase") {
        $dn.RefersTo = "=ID!`$A`$1:`$P`$1"
    }
    elseif ($dn.Name -eq "Measures!_FilterDatabase") {
        $dn.RefersTo = "=Measures!`$A`$1:`$P`$76"
    }
}

# ------------------------------------------------------------------
# 4. Restore per-sheet selections / active cell, then leave NewVars as
#    the active (visible) sheet, matching the saved workbook state.
# ------------------------------------------------------------------
$wsMeasures.Activate() | Out-Null
$wsMeasures.Range("R23").Select() | Out-Null

$wsID = $wb.Worksheets.Item("ID")
$wsID.Activate() | Out-Null
$wsID.Range("O2").Select() | Out-Null

$wsDems = $wb.Worksheets.Item("Dems")
$wsDems.Activate() | Out-Null
$wsDems.Range("O2").Select() | Out-Null

$wsDates = $wb.Worksheets.Item("Dates")
$wsDates.Activate() | Out-Null
$wsDates.Range("O2").Select() | Out-Null

$wsNewVars = $wb.Worksheets.Item("NewVars")
$wsNewVars.Activate() | Out-Null
$wsNewVars.Range("O2").Select() | Out-Null
